$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode subscript-3 character (U+2083), built via code point to avoid PS char/int coercion issues
$sub3 = [string][char]0x2083

function Set-TextCell {
    param($cell, $value)
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value2 = $value
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" ('41.973.98')
Set-TextCell "E2" ('  +0.41%  ')

# Row 3
Set-TextCell "D3" ('2.260.29')
Set-TextCell "E3" ('  -0.36%  ')

# Row 5
Set-TextCell "D5" ('304.75')
Set-TextCell "E5" ('  +0.04%  ')

# Row 6
Set-TextCell "D6" ('95.41')
Set-TextCell "E6" ('  +2.58%  ')

# Row 7
Set-TextCell "D7" ('0.527')
Set-TextCell "E7" ('  -0.44%  ')

# Row 8
Set-TextCell "E8" ('  +0.05%  ')

# Row 9
Set-TextCell "D9" ('0.487')
Set-TextCell "E9" ('  +0.50%  ')

# Row 10
Set-TextCell "D10" ('35.02')
Set-TextCell "E10" ('  +7.57%  ')

# Row 11
Set-TextCell "D11" ('0.0788')
Set-TextCell "E11" ('  -1.16%  ')

# Row 12
Set-TextCell "E12" ('  -0.14%  ')

# Row 13
Set-TextCell "D13" ('6.61')
Set-TextCell "E13" ('  -0.89%  ')

# Row 14
Set-TextCell "D14" ('2.621.32')
Set-TextCell "E14" ('  +0.02%  ')

# Row 15
Set-TextCell "D15" ('14.33')
Set-TextCell "E15" ('  +0.23%  ')

# Row 16
Set-TextCell "D16" ('2.273.53')
Set-TextCell "E16" ('  -0.09%  ')

# Row 17
Set-TextCell "D17" ('0.789')
Set-TextCell "E17" ('  +0.60%  ')

# Row 18
Set-TextCell "D18" ('41.905.79')
Set-TextCell "E18" ('  +0.39%  ')

# Row 19
Set-TextCell "D19" ('12.33')
Set-TextCell "E19" ('  -4.56%  ')

# Row 20
Set-TextCell "D20" ('0.0' + $sub3 + '0901')
Set-TextCell "E20" ('  -0.80%  ')

# Row 21
Set-TextCell "D21" ('5.95')
Set-TextCell "E21" ('  -0.24%  ')

# Row 22
Set-TextCell "D22" ('67.60')
Set-TextCell "E22" ('  -0.18%  ')

# Row 23
Set-TextCell "D23" ('237.07')
Set-TextCell "E23" ('  -2.89%  ')

# Row 24
Set-TextCell "D24" ('2.56')
Set-TextCell "E24" ('  -0.97%  ')

# Row 25
Set-TextCell "B25" ('Dai')
Set-TextCell "C25" ('https://coinranking.com/coin/MoTuySvg7+dai-dai')
Set-TextCell "D25" ('0.999')
Set-TextCell "E25" ('  -0.14%  ')

# Row 26
Set-TextCell "B26" ('ImmutableX')
Set-TextCell "C26" ('https://coinranking.com/coin/Z96jIvLU7+immutablex-imx')
Set-TextCell "D26" ('1.92')
Set-TextCell "E26" ('  -0.84%  ')

# Row 27
Set-TextCell "D27" ('23.60')
Set-TextCell "E27" ('  -1.80%  ')

# Row 28
Set-TextCell "D28" ('36.67')
Set-TextCell "E28" ('  +4.99%  ')

# Row 29
Set-TextCell "D29" ('9.47')
Set-TextCell "E29" ('  -1.37%  ')

# Row 30
Set-TextCell "E30" ('  +1.20%  ')

# Row 31
Set-TextCell "D31" ('159.43')
Set-TextCell "E31" ('  -0.08%  ')

# Row 32
Set-TextCell "D32" ('5.20')
Set-TextCell "E32" ('  -2.13%  ')

# Row 33
Set-TextCell "E33" ('  +0.01%  ')

# Row 34
Set-TextCell "D34" ('3.17')
Set-TextCell "E34" ('  +4.46%  ')

# Row 35
Set-TextCell "D35" ('0.0737')
Set-TextCell "E35" ('  -0.80%  ')

# Row 36
Set-TextCell "D36" ('16.94')
Set-TextCell "E36" ('  +0.29%  ')

# Row 37
Set-TextCell "E37" ('  +0.28%  ')

# Row 38
Set-TextCell "E38" ('  -1.19%  ')

# Row 39
Set-TextCell "D39" ('1.82')
Set-TextCell "E39" ('  +1.23%  ')

# Row 40
Set-TextCell "E40" ('  -1.85%  ')

# Row 41
Set-TextCell "D41" ('4.01')
Set-TextCell "E41" ('  +1.89%  ')

# Row 42
Set-TextCell "D42" ('2.42')
Set-TextCell "E42" ('  +7.90%  ')

# Row 43
Set-TextCell "D43" ('1.979.06')
Set-TextCell "E43" ('  -1.47%  ')

# Row 44
Set-TextCell "D44" ('18.91')
Set-TextCell "E44" ('  -4.56%  ')

# Row 45
Set-TextCell "D45" ('0.0282')
Set-TextCell "E45" ('  -0.09%  ')

# Row 46
Set-TextCell "D46" ('9.89')
Set-TextCell "E46" ('  -3.85%  ')

# Row 47
Set-TextCell "D47" ('2.91')
Set-TextCell "E47" ('  -0.24%  ')

# Row 48
Set-TextCell "D48" ('53.09')
Set-TextCell "E48" ('  -0.79%  ')

# Row 49
Set-TextCell "B49" ('BitcoinSV')
Set-TextCell "C49" ('https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv')
Set-TextCell "D49" ('71.98')
Set-TextCell "E49" ('  -1.75%  ')

# Row 50
Set-TextCell "B50" ('Stacks')
Set-TextCell "C50" ('https://coinranking.com/coin/mMPrMcB7+stacks-stx')
Set-TextCell "D50" ('1.50')
Set-TextCell "E50" ('  +0.14%  ')

# Row 51
Set-TextCell "D51" ('90.64')
Set-TextCell "E51" ('  -1.25%  ')
